$wb = $excel.ActiveWorkbook

# Locate the current last sheet ("two_by_two_PriceinOutput") so the new sheet
# is inserted immediately after it, matching the target sheet order.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Add the new worksheet right after the current last sheet and rename it.
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "two_by_two_PriceinDem"

# Write the "DY" row label first so it is interned into the shared-string
# table ahead of the new header labels below (matches upstream ordering).
$newSheet.Range("A21").Value = "DY"

$arr = New-Object 'object[,]' 28,9
$arr[0,1] = 'benchmark'
$arr[0,2] = 'RA=157'
$arr[0,3] = 'eRA=.5'
$arr[0,4] = 'pr_Ud=2'
$arr[0,5] = 'prU2,eRA.6'
$arr[0,6] = 'prU.5,eRA.6'
$arr[0,7] = 'Itax=0.1'
$arr[0,8] = 'Otax=0.1'
$arr[1,0] = 'X'
$arr[1,1] = 1
$arr[1,2] = 1.0363877164248041
$arr[1,3] = 1.0365519507521972
$arr[1,4] = 1.0595609099982088
$arr[1,5] = 1.0633532225145439
$arr[1,6] = 0.99833470833445848
$arr[1,7] = 0.98248710709607423
$arr[1,8] = 0.94048657015713999
$arr[2,0] = 'Y'
$arr[2,1] = 1
$arr[2,2] = 1.0432700717660239
$arr[2,3] = 1.0430251457731119
$arr[2,4] = 1.0087067867541124
$arr[2,5] = 1.0030495908339756
$arr[2,6] = 1.1000067470042523
$arr[2,7] = 1.122931737696433
$arr[2,8] = 1.1855074833626025
$arr[3,0] = 'U'
$arr[3,1] = 1
$arr[3,2] = 1.0388246300487494
$arr[3,3] = 1.0389899529717939
$arr[3,4] = 1.0621536850625561
$arr[3,5] = 1.0659719353994179
$arr[3,6] = 1.0005254867593421
$arr[3,7] = 0.99843756359281499
$arr[3,8] = 0.99193526198210002
$arr[4,0] = 'PX'
$arr[4,1] = 1
$arr[4,2] = 1.0023513532505028
$arr[4,3] = 1.0023520310092886
$arr[4,4] = 1.0024470278582813
$arr[4,5] = 1.0024626933265708
$arr[4,6] = 1.0021944327932055
$arr[4,7] = 1.0162347743614695
$arr[4,8] = 1.054704334391892
$arr[5,0] = 'PY'
$arr[5,1] = 1
$arr[5,2] = 0.99573893485424703
$arr[5,3] = 0.99573771069529127
$arr[5,4] = 0.99556615181885155
$arr[5,5] = 0.99553786527941868
$arr[5,6] = 0.99602242514969408
$arr[5,7] = 0.97114380037066672
$arr[5,8] = 0.90770374205596749
$arr[6,0] = 'PU'
$arr[6,1] = 1
$arr[6,2] = 1
$arr[6,3] = 1
$arr[6,4] = 1
$arr[6,5] = 1
$arr[6,6] = 1
$arr[6,7] = 1
$arr[6,8] = 1
$arr[7,0] = 'PL'
$arr[7,1] = 1
$arr[7,2] = 0.94438602708152408
$arr[7,3] = 0.94437046939986447
$arr[7,4] = 0.94219250114219932
$arr[7,5] = 0.94183384552627059
$arr[7,6] = 0.9479952509302686
$arr[7,7] = 0.89903690333116493
$arr[7,8] = 0.84467427236056358
$arr[8,0] = 'PK'
$arr[8,1] = 1
$arr[8,2] = 1.0388246300416317
$arr[8,3] = 1.0388360220635884
$arr[8,4] = 1.0404339284158413
$arr[8,5] = 1.0406976527852307
$arr[8,6] = 1.0361901997164102
$arr[8,7] = 1.0329715504842176
$arr[8,8] = 0.96149717694622527
$arr[9,0] = 'SX'
$arr[9,1] = 80
$arr[9,2] = 80
$arr[9,3] = 80
$arr[9,4] = 80
$arr[9,5] = 80
$arr[9,6] = 80
$arr[9,7] = 80
$arr[9,8] = 80
$arr[10,0] = 'SY'
$arr[10,1] = 54
$arr[10,2] = 54
$arr[10,3] = 54
$arr[10,4] = 54
$arr[10,5] = 54
$arr[10,6] = 54
$arr[10,7] = 54
$arr[10,8] = 54.000000000000007
$arr[11,0] = 'SU'
$arr[11,1] = 124
$arr[11,2] = 124
$arr[11,3] = 123.99999999999999
$arr[11,4] = 124
$arr[11,5] = 124
$arr[11,6] = 124
$arr[11,7] = 124
$arr[11,8] = 124
$arr[12,0] = 'DXL'
$arr[12,1] = 30
$arr[12,2] = 31.841365430850797
$arr[12,3] = 31.841911518828791
$arr[12,4] = 31.918541910798012
$arr[12,5] = 31.931195659031214
$arr[12,6] = 31.715172575173078
$arr[12,7] = 30.827982414866039
$arr[12,8] = 30.648734246825121
$arr[13,0] = 'DXK'
$arr[13,1] = 50
$arr[13,2] = 48.244493065345502
$arr[13,3] = 48.243996628540472
$arr[13,4] = 48.174468386692752
$arr[13,5] = 48.16301308279445
$arr[13,6] = 48.359578823666638
$arr[13,7] = 49.189872358276297
$arr[13,8] = 49.36228226730212
$arr[14,0] = 'DYL'
$arr[14,1] = 24
$arr[14,2] = 25.305048730086334
$arr[14,3] = 25.305434497175867
$arr[14,4] = 25.359560402662385
$arr[14,5] = 25.368496662333634
$arr[14,6] = 25.215883919305377
$arr[14,7] = 25.924910448641455
$arr[14,8] = 25.790876462309562
$arr[15,0] = 'DYK'
$arr[15,1] = 30
$arr[15,2] = 28.755737186305545
$arr[15,3] = 28.755386493763073
$arr[15,4] = 28.706277005056187
$arr[15,5] = 28.698187104056096
$arr[15,6] = 28.837054010514596
$arr[15,7] = 28.204372131496182
$arr[15,8] = 28.321572766511753
$arr[16,0] = 'DUX'
$arr[16,1] = 80
$arr[16,2] = 79.812333011330708
$arr[16,3] = 79.812279044731554
$arr[16,4] = 79.804715637612503
$arr[16,5] = 79.803468530612463
$arr[16,6] = 79.824829775826657
$arr[16,7] = 78.721966634399251
$arr[16,8] = 75.850641162032716
$arr[17,0] = 'DUY'
$arr[17,1] = 44
$arr[17,2] = 44.188289178863421
$arr[17,3] = 44.18834350389114
$arr[17,4] = 44.195958168740567
$arr[17,5] = 44.197213922797879
$arr[17,6] = 44.175712201847752
$arr[17,7] = 45.307399360636374
$arr[17,8] = 48.473965635901294
$arr[18,0] = 'RA'
$arr[18,1] = 134
$arr[18,2] = 139.20250041197306
$arr[18,3] = 139.202487647439
$arr[18,4] = 139.20094884111393
$arr[18,5] = 139.20074264707893
$arr[18,6] = 139.20613388257078
$arr[18,7] = 138.76352763367868
$arr[18,8] = 137.46365821827374
$arr[19,0] = 'DU'
$arr[19,1] = 124
$arr[19,2] = 128.81425411257209
$arr[19,3] = 128.83475414329422
$arr[19,4] = 131.70705694774205
$arr[19,5] = 132.18051998952777
$arr[19,6] = 124.06516035720092
$arr[19,7] = 123.80625788538022
$arr[19,8] = 122.99997248577242
$arr[20,0] = 'DY'
$arr[20,1] = 10
$arr[20,2] = 10.432700716801408
$arr[20,3] = 10.412112941775911
$arr[20,4] = 7.5272666509210593
$arr[20,5] = 7.0516882404877128
$arr[20,6] = 15.201438384376063
$arr[20,7] = 15.401704405248296
$arr[20,8] = 15.934368299221509
$arr[21,0] = 'CWI'
$arr[21,1] = 1
$arr[21,2] = 1.0391557242667655
$arr[21,3] = 1.0391553971951093
$arr[21,4] = 1.0454956618531788
$arr[21,5] = 1.046614299839608
$arr[21,6] = 1.0534806670456223
$arr[21,7] = 1.0530112713727526
$arr[21,8] = 1.0506793942727906
$arr[22,0] = 'PX/PX'
$arr[22,1] = 1
$arr[22,2] = 1
$arr[22,3] = 1
$arr[22,4] = 1
$arr[22,5] = 1
$arr[22,6] = 1
$arr[22,7] = 1
$arr[22,8] = 1
$arr[23,0] = 'PY/PX'
$arr[23,1] = 1
$arr[23,2] = 0.99340309326184617
$arr[23,3] = 0.99340120026759737
$arr[23,4] = 0.99313592055419553
$arr[23,5] = 0.99309218378573993
$arr[23,6] = 0.9938415067560199
$arr[23,7] = 0.95562937312479324
$arr[23,8] = 0.86062388525151912
$arr[24,0] = 'PU/PX'
$arr[24,1] = 1
$arr[24,2] = 0.99765416264179452
$arr[24,3] = 0.997653488059559
$arr[24,4] = 0.99755894547015678
$arr[24,5] = 0.9975433566326557
$arr[24,6] = 0.99781037219784841
$arr[24,7] = 0.98402458292999251
$arr[24,8] = 0.94813301452540932
$arr[25,0] = 'PL/PX'
$arr[25,1] = 1
$arr[25,2] = 0.94217065105862896
$arr[25,3] = 0.94215449281721775
$arr[25,4] = 0.93989255786930181
$arr[25,5] = 0.93952009565651817
$arr[25,6] = 0.94591949417252408
$arr[25,7] = 0.88467441383912149
$arr[25,8] = 0.80086356414527782
$arr[26,0] = 'PK/PX'
$arr[26,1] = 1
$arr[26,2] = 1.0363877164158561
$arr[26,3] = 1.0363983809336559
$arr[26,4] = 1.0378941724618791
$arr[26,5] = 1.0381410297991052
$arr[26,6] = 1.0339213288467941
$arr[26,7] = 1.0164693991437799
$arr[26,8] = 0.91162721683569548
$arr[27,0] = 'RA/PX'
$arr[27,1] = 134
$arr[27,2] = 138.87595398615105
$arr[27,3] = 138.87584734803519
$arr[27,4] = 138.86115173438685
$arr[27,5] = 138.85877606592558
$arr[27,6] = 138.90132426159147
$arr[27,7] = 136.54672240562513
$arr[27,8] = 130.33383265418243


$newSheet.Range("A1").Resize(28, 9).Value = $arr
